# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
# These two sheets carry duplicate rows of exhibition data; the "全部类型"
# sheet has one extra row inserted around row 33, so matching rows there
# are offset by +1 from row 33 onward.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row (in "展览") -> new value for column F
$updatesExhibit = @{
    2  = 3161
    3  = 565
    6  = 85
    10 = 16462
    11 = 283
    14 = 6391
    15 = 641
    16 = 130
    22 = 39
    29 = 901
    31 = 5066
    33 = 11378
    38 = 3847
}

# Row (in "全部类型") -> new value for column F
$updatesAll = @{
    2  = 3161
    3  = 565
    6  = 85
    10 = 16462
    11 = 283
    14 = 6391
    15 = 641
    16 = 130
    22 = 39
    29 = 901
    31 = 5066
    34 = 11378
    39 = 3847
}

foreach ($row in $updatesExhibit.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $updatesExhibit[$row]
}

foreach ($row in $updatesAll.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $updatesAll[$row]
}
